$d = $word.ActiveDocument

# 1. Update the "Nombre" cell text to mention the extended use case.
$d.Content.Find.Execute("CU-06 Eliminar producto ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "CU-06 Eliminar producto extendido de CU-04 Ver producto", 2)

# 2. Add a new "Extiende" / "CU-04" row at the end of the table
#    (mirrors the row just above it, e.g. "Postcondiciones").
$t = $d.Tables.Item(1)
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "Extiende"
$newRow.Cells.Item(2).Range.Text = "CU-04"
